$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.430.17'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.573.35'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.13'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +2.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.89'
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3427'
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07666'
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.24'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.943'
$ws.Range("D16").Value = '1.570.78'
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.33'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06767'
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("E21").Value = '  +2.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.227'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.04'
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.429'
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("D25").Value = '22.431.38'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.733'
$ws.Range("E26").Value = '  -8.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.31'
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '146.07'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.034'
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.46'
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("D31").Value = '1.746.36'
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.205'
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.018'
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("E35").Value = '  -3.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08553'
$ws.Range("E36").Value = '  +1.44%  '
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2318'
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("E39").Value = '  +1.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.341'
$ws.Range("E40").Value = '  +7.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.477'
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.57'
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6453'
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.12'
$ws.Range("E44").Value = '  -3.28%  '
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.799'
$ws.Range("E46").Value = '  +0.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6016'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.302'
$ws.Range("E48").Value = '  +8.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.088'
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.44'
$ws.Range("E50").Value = '  +3.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07329'
